# Add a "Saves" column to the earnings breakdown table.
#
# Current layout:  A Member | B Respect | C War_Hits | D Save_Score | E Total_Score | F..K (shares/earnings)
# New layout:      A Member | B Respect | C War_Hits | D Total_Score | E Saves | F Save_Score | G..L (shares/earnings)
#
# Each "save" is worth a constant 11.86 points, so Saves = Save_Score / 11.86.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at E; this shifts the old E:K block to F:L
# (old Total_Score -> F, old Respect_Share.. -> G..L automatically).
$ws.Columns("E").Insert()

# Fix up the three header cells that actually changed meaning/position.
$ws.Range("D1").Value = "Total_Score"
$ws.Range("E1").Value = "Saves"
$ws.Range("F1").Value = "Save_Score"

# Find the last populated data row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# D currently still holds the old Save_Score values, F currently holds the
# old Total_Score values (shifted there by the column insert). Read both,
# then recompute D/E/F for the new layout in one bulk read/write.
$rowCount = $lastRow - 1
$srcRange = $ws.Range("D2:F" + $lastRow)
$src = $srcRange.Value2

$result = New-Object 'object[,]' $rowCount,3
for ($i = 1; $i -le $rowCount; $i++) {
    $saveScore = $src[$i,1]
    $totalScore = $src[$i,3]
    if ($null -eq $saveScore) { $saveScore = 0 }
    if ($null -eq $totalScore) { $totalScore = 0 }
    $saves = [math]::Round($saveScore / 11.86)

    $result[$i-1,0] = $totalScore
    $result[$i-1,1] = $saves
    $result[$i-1,2] = $saveScore
}

$ws.Range("D2:F" + $lastRow).Value2 = $result

"done"
